$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new 2020 data point in column Q, mirroring the formatting
# already used for the 2019 column (P).
$ws.Range("P4").Copy()
$ws.Range("Q4").PasteSpecial(-4122)
$ws.Range("Q4").Value = 2020

$ws.Range("P5").Copy()
$ws.Range("Q5").PasteSpecial(-4122)
$ws.Range("Q5").Value = 90.6

$excel.CutCopyMode = 0

# Leave the selection where the editor ended up after adding the data.
$ws.Range("P12").Select()
